# Generate Report for Handback
# Updates timestamps / status text produced by the handback report generator.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" column (G) on the Overview sheet.
$wsOverview.Range("G2").Value = "2016-08-17 14:14:01"
$wsOverview.Range("G3").Value = "2016-08-17 14:14:01"

# zh-cn sheet: Priority (E) changes from "ht" to "mt"; Correspond Handoff
# Datetime (H) and Correspond Handback DateTime (K) move forward.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-17 14:13:55"
$wsZhCn.Range("H3").Value = "2016-08-17 14:13:55"
$wsZhCn.Range("K2").Value = "2016-08-17 14:14:29"
$wsZhCn.Range("K3").Value = "2016-08-17 14:14:29"

# de-de sheet: Priority (E) changes from "ht" to "mt"; Correspond Handback
# DateTime (K) moves forward. Correspond Handoff Datetime (H) reuses the same
# shared value as the Overview sheet's "Latest HO Xliff Generate Date".
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-17 14:14:01"
$wsDeDe.Range("H3").Value = "2016-08-17 14:14:01"
$wsDeDe.Range("K2").Value = "2016-08-17 14:14:37"
$wsDeDe.Range("K3").Value = "2016-08-17 14:14:37"
